$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "/media/Data/common/guideseq_gnt_dev/test_dataset/"
$ws.Range("O1").Value = "path_to_files"
$ws.Range("O3").Value = "/media/Data/common/guideseq_gnt_dev/test_dataset/"

$ws.Range("O2").Select()
